$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the rate-conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$oldLine1 = [string]::Concat([char]0x2705, " 1000 Bs = 3.25 = 12534.15 pesos")
$newLine1 = [string]::Concat([char]0x2705, " 1000 Bs = 3.25 = 12494.84 pesos")
$oldLine2 = [string]::Concat([char]0x2705, " 12534.15 pesos = 3.23 = 968.89 Bs")
$newLine2 = [string]::Concat([char]0x2705, " 12494.84 pesos = 3.24 = 980.71 Bs")

$cellA1 = $wsHoja1.Range("A1")
$text = $cellA1.Value()
$text = $text.Replace($oldLine1, $newLine1)
$text = $text.Replace($oldLine2, $newLine2)
$cellA1.Value = $text

# --- Sheet "tasas": update N10/O10 and N12/O12 values ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 308
$wsTasas.Range("O10").Value = 3848.41

$wsTasas.Range("N12").Value = 3858
$wsTasas.Range("O12").Value = 302.81
